# Update gh-pages output data (statistics refresh) across the four sheets:
# 展览 (Exhibition), 演出 (Performance), 本地生活 (Local Life), 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet: 展览 (Exhibition) - "F" column (想去人数 / interest count) refresh
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 662
$ws.Range("F3").Value = 10805
$ws.Range("F4").Value = 7183
$ws.Range("F5").Value = 147
$ws.Range("F6").Value = 13257
$ws.Range("F7").Value = 13405
$ws.Range("F8").Value = 1343
$ws.Range("F9").Value = 1322
$ws.Range("F10").Value = 5626
$ws.Range("F11").Value = 942
$ws.Range("F12").Value = 554
$ws.Range("F16").Value = 1468
$ws.Range("F17").Value = 386
$ws.Range("F18").Value = 2070
$ws.Range("F19").Value = 1096
$ws.Range("F20").Value = 1666
$ws.Range("F23").Value = 2201
$ws.Range("F24").Value = 530
$ws.Range("F25").Value = 756
$ws.Range("F26").Value = 3111
$ws.Range("F28").Value = 2173
$ws.Range("F29").Value = 30
$ws.Range("F31").Value = 1722
$ws.Range("F32").Value = 1032
$ws.Range("F33").Value = 1205
$ws.Range("F34").Value = 76
$ws.Range("F35").Value = 117
$ws.Range("F36").Value = 4400
$ws.Range("F37").Value = 4542
$ws.Range("F38").Value = 283
$ws.Range("F39").Value = 142
$ws.Range("F40").Value = 629
$ws.Range("F42").Value = 3183
$ws.Range("F46").Value = 58
$ws.Range("F47").Value = 47
$ws.Range("F48").Value = 4341
$ws.Range("F49").Value = 223

# ------------------------------------------------------------------
# Sheet: 演出 (Performance)
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 94
$ws.Range("F17").Value = 16
$ws.Range("F19").Value = 15
$ws.Range("F21").Value = 82
$ws.Range("F26").Value = 67
$ws.Range("F27").Value = 9
# These two shows became unavailable for ticket purchase (numeric lowest-price -> text)
$ws.Range("G2").Value = "不可售"
$ws.Range("G3").Value = "不可售"

# ------------------------------------------------------------------
# Sheet: 本地生活 (Local Life)
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
# Row 2 re-opened for sale: time window extended, interest count refreshed,
# and lowest price switched from sold-out text back to a numeric price.
$ws.Range("E2").Value = "2024.06.08 00:00-09.08 23:59"
$ws.Range("F2").Value = 7091
$ws.Range("G2").Value = 10
$ws.Range("F3").Value = 148
$ws.Range("F4").Value = 428

# ------------------------------------------------------------------
# Sheet: 全部类型 (All Types)
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 662
$ws.Range("F3").Value = 10805
$ws.Range("F4").Value = 7183
$ws.Range("F5").Value = 147
$ws.Range("F6").Value = 148
$ws.Range("F7").Value = 428
$ws.Range("F8").Value = 13257
$ws.Range("F9").Value = 13405
$ws.Range("F11").Value = 1343
$ws.Range("F12").Value = 1322
$ws.Range("F13").Value = 5626
$ws.Range("F15").Value = 94
$ws.Range("F18").Value = 1468
$ws.Range("F19").Value = 386
$ws.Range("F20").Value = 2070
$ws.Range("F21").Value = 1096
$ws.Range("F22").Value = 1666
$ws.Range("F24").Value = 530
$ws.Range("F25").Value = 756
$ws.Range("F26").Value = 3111
$ws.Range("F29").Value = 2173
$ws.Range("F30").Value = 30
$ws.Range("F32").Value = 1722
$ws.Range("F33").Value = 16
$ws.Range("F34").Value = 1032
$ws.Range("F35").Value = 1205
$ws.Range("F36").Value = 117
$ws.Range("F38").Value = 4400
$ws.Range("F39").Value = 4542
$ws.Range("F40").Value = 283
$ws.Range("F41").Value = 142
$ws.Range("F42").Value = 629
$ws.Range("F44").Value = 3183
$ws.Range("F47").Value = 58
$ws.Range("F48").Value = 47
$ws.Range("F49").Value = 9
